$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Caris LeVert", "SG,SF", "Cleveland Cavaliers"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Ja Morant", "PG", "Memphis Grizzlies")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
